$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "length (optional column, for checking only)"
$ws.Range("H1").Value = "answer (optional column, for checking only)"

$ws.Range("H2").Select()
